$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.ClearFormats()
}

Set-TextValue 'D2' '70.408.18'
Set-TextValue 'E2' '  +0.73%  '
Set-TextValue 'D3' '3.572.85'
Set-TextValue 'E3' '  -0.11%  '
Set-TextValue 'E4' '  +0.01%  '
Set-TextValue 'D5' '590.84'
Set-TextValue 'E5' '  +2.83%  '
Set-TextValue 'D6' '187.29'
Set-TextValue 'E6' '  +0.31%  '
Set-TextValue 'D7' '3.558.96'
Set-TextValue 'E7' '  -0.43%  '
Set-TextValue 'D8' '0.621'
Set-TextValue 'E8' '  -0.31%  '
Set-TextValue 'E9' '  +0.01%  '
Set-TextValue 'E10' '  +8.20%  '
Set-TextValue 'E11' '  -0.45%  '
Set-TextValue 'D12' '54.85'
Set-TextValue 'E12' '  -0.75%  '
Set-TextValue 'D13' '0.0000308'
Set-TextValue 'E13' '  +0.71%  '
Set-TextValue 'D14' '9.58'
Set-TextValue 'E14' '  +0.06%  '
Set-TextValue 'D15' '4.139.97'
Set-TextValue 'E15' '  -0.31%  '
Set-TextValue 'D16' '19.49'
Set-TextValue 'E16' '  -0.82%  '
Set-TextValue 'D17' '70.379.19'
Set-TextValue 'E17' '  +0.73%  '
Set-TextValue 'D18' '3.592.81'
Set-TextValue 'E18' '  +0.46%  '
Set-TextValue 'E19' '  -0.24%  '
Set-TextValue 'E20' '  -0.81%  '
Set-TextValue 'D21' '548.94'
Set-TextValue 'E21' '  +11.93%  '
Set-TextValue 'E22' '  -0.81%  '
Set-TextValue 'D23' '18.02'
Set-TextValue 'E23' '  -7.87%  '
Set-TextValue 'D24' '4.67'
Set-TextValue 'E24' '  +8.35%  '
Set-TextValue 'D25' '4.92'
Set-TextValue 'E25' '  -0.36%  '
Set-TextValue 'D26' '96.11'
Set-TextValue 'E26' '  +0.11%  '
Set-TextValue 'D27' '11.49'
Set-TextValue 'E27' '  +3.81%  '
Set-TextValue 'D28' '3.01'
Set-TextValue 'E28' '  +1.84%  '
Set-TextValue 'D29' '9.17'
Set-TextValue 'E29' '  -1.19%  '
Set-TextValue 'D31' '7.36'
Set-TextValue 'E31' '  -2.27%  '
Set-TextValue 'D32' '12.58'
Set-TextValue 'E32' '  +4.47%  '
Set-TextValue 'D33' '65.21'
Set-TextValue 'E33' '  -2.54%  '
Set-TextValue 'E34' '  -0.89%  '
Set-TextValue 'D35' '551.24'
Set-TextValue 'E35' '  -3.34%  '
Set-TextValue 'D36' '3.24'
Set-TextValue 'E36' '  +5.34%  '
Set-TextValue 'D37' '0.419'
Set-TextValue 'E37' '  +6.79%  '
Set-TextValue 'E38' '  +0.23%  '
Set-TextValue 'E39' '  +0.05%  '
Set-TextValue 'D40' '0.0₃0771'
Set-TextValue 'E40' '  -3.61%  '
Set-TextValue 'E41' '  -0.86%  '
Set-TextValue 'D42' '3.373.35'
Set-TextValue 'E42' '  +3.77%  '
Set-TextValue 'E43' '  -3.88%  '
Set-TextValue 'D44' '3.09'
Set-TextValue 'E44' '  -6.29%  '
Set-TextValue 'D45' '3.59'
Set-TextValue 'E45' '  +3.76%  '
Set-TextValue 'B46' 'VeChain'
Set-TextValue 'C46' 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 'D46' '0.0449'
Set-TextValue 'E46' '  +3.24%  '
Set-TextValue 'B47' 'ThetaToken'
Set-TextValue 'C47' 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
Set-TextValue 'D47' '3.00'
Set-TextValue 'E47' '  +0.20%  '
Set-TextValue 'D48' '9.23'
Set-TextValue 'E48' '  -4.32%  '
Set-TextValue 'E49' '  +0.19%  '
Set-TextValue 'D50' '0.998'
Set-TextValue 'E50' '  -0.01%  '
Set-TextValue 'D51' '1.47'
Set-TextValue 'E51' '  +21.48%  '
